$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.316.83'
$ws.Range("E2").Value = '  -0.68%  '

$ws.Range("D3").Value = '1.859.63'
$ws.Range("E3").Value = '  -2.01%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '235.61'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.48%  '

$ws.Range("E6").Value = '  -0.08%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4782'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.16%  '

$ws.Range("E8").Value = '  -3.75%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06464'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.97%  '

$ws.Range("D10").Value = '1.859.42'
$ws.Range("E10").Value = '  -2.16%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07398'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.11%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '16.23'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -3.75%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.014'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.24%  '

$ws.Range("E14").Value = '  -1.39%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6448'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.81%  '

$ws.Range("D16").Value = '30.275.17'
$ws.Range("E16").Value = '  -0.64%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.002'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.19%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.15'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.70%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007531'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.68%  '

$ws.Range("B20").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C20").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D20").Value = '2.093.31'
$ws.Range("E20").Value = '  -2.42%  '

$ws.Range("B21").Value = 'BinanceUSD'
$ws.Range("C21").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.001'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.02%  '

$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.270'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.14%  '

$ws.Range("B23").Value = 'BitcoinCash'
$ws.Range("C23").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '218.54'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +12.08%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.064'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.29%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.216'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.26%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '164.50'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.61%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.44'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.39%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.925'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.11%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.430'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.02%  '

$ws.Range("B30").Value = 'Stellar'
$ws.Range("C30").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.09261'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.00%  '

$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.265'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.22%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.955'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.89%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04977'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.75%  '

$ws.Range("B34").Value = 'ImmutableX'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7294'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.08%  '

$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.143'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.24%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.688'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.47%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01820'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.43%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.607'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.68%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.8997'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.74%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.036'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.31%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.891'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.07%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '105.98'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.82%  '

$ws.Range("B43").Value = 'PaxDollar'
$ws.Range("C43").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.001'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.65%  '

$ws.Range("B44").Value = 'TheSandbox'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.4243'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.05%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '7.312'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.56%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.1292'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.80%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '63.02'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -6.52%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.485'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +7.05%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.805'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.16%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '33.74'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.85%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05637'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.24%  '
